# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Membrillo" (Terminal La Palmera de La Serena)
# above the former row 28, pushing the existing rows 28-47 down to rows 31-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 28 (new rows inherit the date
# number format from the row above, matching column D's existing formatting).
$ws.Rows("28:30").Insert()

# New data for rows 28, 29 and 30 (columns A through T).
$newRows = @(
    ,@(8, "Terminal La Palmera de La Serena", "Coquimbo", 44679, 4, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Especial", 16, 320000, 330000, 325000, "`$/bins (450 kilos)", "Provincia de Cachapoal", 722, 450)
    ,@(8, "Terminal La Palmera de La Serena", "Coquimbo", 44679, 4, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Primera", 16, 280000, 290000, 285000, "`$/bins (450 kilos)", "Provincia de Cachapoal", 633, 450)
    ,@(8, "Terminal La Palmera de La Serena", "Coquimbo", 44679, 4, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Segunda", 16, 250000, 260000, 255000, "`$/bins (450 kilos)", "Provincia de Cachapoal", 567, 450)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowData = $newRows[$i]
    $targetRow = 28 + $i
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($targetRow, $col).Value = $rowData[$j]
    }
}
